$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2407
$ws.Range("I70").Value = 2552.8572
$ws.Range("J70").Value = 2066.6667
$ws.Range("K70").Value = 7658.571599999999
$ws.Range("L70").Value = 6200.000100000001
$ws.Range("M70").Value = -7388.571599999999
$ws.Range("N70").Value = -6740.000100000001
$ws.Range("H73").Value = 2407
$ws.Range("I73").Value = 2552.8572
$ws.Range("J73").Value = 2066.6667
$ws.Range("K73").Value = 7658.571599999999
$ws.Range("L73").Value = 6200.000100000001
$ws.Range("M73").Value = -6722.571599999999
$ws.Range("N73").Value = -8072.000100000001
$ws.Range("H97").Value = 1985.625
$ws.Range("J97").Value = 2158.2856
$ws.Range("L97").Value = 6474.8568
$ws.Range("N97").Value = -7466.8568
$ws.Range("H99").Value = 111111690
$ws.Range("I99").Value = 626.4286
$ws.Range("K99").Value = 1879.2858
$ws.Range("M99").Value = -381.2857999999999
$ws.Range("H100").Value = 3862.7273
$ws.Range("I100").Value = 2886.25
$ws.Range("K100").Value = 2886.25
$ws.Range("M100").Value = -2345.25
$ws.Range("H101").Value = 55557720
$ws.Range("I101").Value = 2949.5
$ws.Range("J101").Value = 166667260
$ws.Range("K101").Value = 8848.5
$ws.Range("L101").Value = 500001780
$ws.Range("M101").Value = -7226.5
$ws.Range("N101").Value = -500005024
$ws.Range("H114").Value = 39561
$ws.Range("J114").Value = 39561
$ws.Range("L114").Value = 39561
$ws.Range("N114").Value = -48239
$ws.Range("H129").Value = 1906
$ws.Range("I129").Value = 784.4167
$ws.Range("J129").Value = 4149.1665
$ws.Range("K129").Value = 2353.2501
$ws.Range("L129").Value = 12447.4995
$ws.Range("M129").Value = 2646.7499
$ws.Range("N129").Value = -22447.4995
$ws.Range("H138").Value = 2893.6206
$ws.Range("I138").Value = 9518.333000000001
$ws.Range("K138").Value = 28554.999
$ws.Range("M138").Value = -23414.999
$ws.Range("H141").Value = 7450.143
$ws.Range("I141").Value = 7247.5
$ws.Range("K141").Value = 21742.5
$ws.Range("M141").Value = -16562.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17278.578
$ws.Range("I32").Value = 3149.3281
$ws.Range("J32").Value = 146460.28
$ws.Range("K32").Value = 3149.3281
$ws.Range("L32").Value = 146460.28
$ws.Range("M32").Value = -2862.3281
$ws.Range("N32").Value = -147034.28
$ws.Range("H61").Value = 2389.8157
$ws.Range("I61").Value = 2070.1538
$ws.Range("K61").Value = 2070.1538
$ws.Range("M61").Value = -1858.1538
$ws.Range("H132").Value = 2158.6511
$ws.Range("I132").Value = 1156.5143
$ws.Range("K132").Value = 3469.5429
$ws.Range("M132").Value = -939.5429000000004
$ws.Range("H136").Value = 2389.8157
$ws.Range("I136").Value = 2070.1538
$ws.Range("K136").Value = 6210.4614
$ws.Range("M136").Value = -3660.4614
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 34999.832
$ws.Range("J35").Value = 34999.832
$ws.Range("L35").Value = 34999.832
$ws.Range("N35").Value = -35619.832
$ws.Range("H76").Value = 19078.5
$ws.Range("I76").Value = 15500
$ws.Range("J76").Value = 22657
$ws.Range("K76").Value = 15500
$ws.Range("L76").Value = 22657
$ws.Range("M76").Value = -15185
$ws.Range("N76").Value = -23287
$ws.Range("H79").Value = 19078.5
$ws.Range("I79").Value = 15500
$ws.Range("J79").Value = 22657
$ws.Range("K79").Value = 15500
$ws.Range("L79").Value = 22657
$ws.Range("M79").Value = -14408
$ws.Range("N79").Value = -24841
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 22080
$ws.Range("I62").Value = 25694.8
$ws.Range("K62").Value = 25694.8
$ws.Range("M62").Value = -25070.8
$ws.Range("H65").Value = 22080
$ws.Range("I65").Value = 25694.8
$ws.Range("K65").Value = 128474
$ws.Range("M65").Value = -125354
$ws.Range("H93").Value = 24315.166
$ws.Range("J93").Value = 28632
$ws.Range("L93").Value = 28632
$ws.Range("N93").Value = -32376
$ws.Range("H134").Value = 3013.8572
$ws.Range("I134").Value = 3151.3333
$ws.Range("K134").Value = 9453.999899999999
$ws.Range("M134").Value = -6918.999899999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 588.3
$ws.Range("J33").Value = 814.2857
$ws.Range("L33").Value = 4885.7142
$ws.Range("N33").Value = -5451.7142
$ws.Range("H86").Value = 487
$ws.Range("I86").Value = 487
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1461
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -275
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 487
$ws.Range("I89").Value = 487
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 4383
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 1545
$ws.Range("N89").ClearContents()
$ws.Range("H109").Value = 1443.6666
$ws.Range("I109").Value = 1443.6666
$ws.Range("K109").Value = 4330.9998
$ws.Range("M109").Value = -3290.9998
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 609.13635
$ws.Range("I2").Value = 326.86667
$ws.Range("J2").Value = 1214
$ws.Range("K2").Value = 326.86667
$ws.Range("L2").Value = 1214
$ws.Range("M2").Value = -213.86667
$ws.Range("N2").Value = -1440
$ws.Range("H132").Value = 5016.6484
$ws.Range("I132").Value = 4101.8213
$ws.Range("J132").Value = 7862.778
$ws.Range("K132").Value = 12305.4639
$ws.Range("L132").Value = 23588.334
$ws.Range("M132").Value = -9775.463899999999
$ws.Range("N132").Value = -28648.334
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7537.6924
$ws.Range("I40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("M40").Value = -4864
$ws.Range("H63").Value = 23085
$ws.Range("J63").Value = 23085
$ws.Range("L63").Value = 23085
$ws.Range("N63").Value = -24583
$ws.Range("H66").Value = 23085
$ws.Range("J66").Value = 23085
$ws.Range("L66").Value = 69255
$ws.Range("N66").Value = -76743
$ws.Range("H93").Value = 31757.166
$ws.Range("J93").Value = 86909.25
$ws.Range("L93").Value = 86909.25
$ws.Range("N93").Value = -89405.25
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
